$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Move the "Category_Attribute_Definition" mini-table from E11:E13 down to E17:E19
#    to make room for the four new attribute rows appended to the Category /
#    Attribute_Definition tables above it.
$ws.Range("E11:E13").Cut($ws.Range("E17:E19"))
$ws.Range("E11:E13").Clear()
$ws.Rows(11).AutoFit()
$ws.Rows(17).RowHeight = 19

# 2. Add the new common audit columns to the Category table (D column) and the
#    Attribute_Definition table (G column), rows 8-12, re-using the existing
#    "bordered, no fill" look used elsewhere in these two tables.
$dgRanges = @("D8:D12", "G8:G12")
$dgValues = @("is_active", "is_deleted", "created_by", "created_at", "updated_at")

foreach ($rng in $dgRanges) {
    $cells = $ws.Range($rng).Cells
    for ($i = 1; $i -le 5; $i++) {
        $cells.Item($i).Value = $dgValues[$i - 1]
    }
    $ws.Range("D6").Copy()
    $ws.Range($rng).PasteSpecial(-4122)
    $ws.Range($rng).Interior.ColorIndex = -4142
}

# 3. Add the same five new audit columns to the Item table (I column), rows 16-20.
$iValues = @("is_active", "is_deleted", "created_by", "created_at", "updated_at")
$iCells = $ws.Range("I16:I20").Cells
for ($i = 1; $i -le 5; $i++) {
    $iCells.Item($i).Value = $iValues[$i - 1]
}
$ws.Range("D6").Copy()
$ws.Range("I16:I20").PasteSpecial(-4122)
$ws.Range("I16:I20").Interior.ColorIndex = -4142

# 4. Widen column D now that it holds the new, longer field names.
$ws.Columns("D").ColumnWidth = 19.1640625

# 5. Restore the selection to where the user last left off.
$ws.Range("E22").Select()
